$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Update Runmode column (C2:C18) from "N" to "Y"
$ws.Range("C2:C18").Value = "Y"

# C18 previously had a distinct "fill" style; align it with the rest of the column (C2:C17)
$ws.Range("C18").Interior.Pattern = -4142

# Update the active selection to reflect the edited range
$ws.Range("C2:C18").Select()
